# V. 115 "Tiempo para mi"
# Adds a new movie row to the "Películas" sheet / "Tabla24" table, keeping
# the table sorted descending by the "Puntuación total" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Películas")

# New entry fits (by score, descending) right above the former row 132
# ("Maldita suerte"), so insert a fresh worksheet row there - this shifts
# the existing rows 132-135 down to 133-136 and adjusts formula refs.
$ws.Range("B132").EntireRow.Insert()

# Fill in the new row's data.
$ws.Range("B132").Value = "Tiempo para mi"
$ws.Range("D132").Value = 5
$ws.Range("E132").Value = 2
$ws.Range("F132").Value = 2
$ws.Range("G132").Value = 3
$ws.Range("H132").Value = 5.0999999999999996
$ws.Range("I132").Value = 4.0999999999999996
$ws.Range("C132").Formula = "=AVERAGE(D132,E132,E132,F132,G132,H132,H132,I132)"

# Grow the table ("Tabla24") so it covers the newly inserted row.
$tbl = $ws.ListObjects.Item("Tabla24")
$tbl.Resize($ws.Range("B2:I136"))

# The "newest entry" highlight moves from the previous newest row (81) to
# the row just added (132).
$ws.Range("B81").HorizontalAlignment = -4131
$ws.Range("B132").HorizontalAlignment = -4131

# Move the window selection to reflect where the edit was made.
$ws.Range("C136").Select()

$wb.Save()
